# Arquivos traduzidos por juminho
# Marks a batch of previously-untranslated rows as translated ("SIM") and
# fills in the translator column ("juminho") for another batch of rows
# that were already marked untranslated but had no translator set.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows whose "Traduzido" (column C) flips from "Não" to "SIM".
# (Column D on these rows already says "juminho".)
$rowsMarkSim = @(339,340,341,342,343,345,346,347,348,349,350,352,353,354,356,357,359,360,361,362,363,364,365,366,367,368,369,370,371,372)

foreach ($r in $rowsMarkSim) {
    $ws.Cells.Item($r, 3).Value = "SIM"
}

# Rows that stay "Não" in column C but now get "juminho" recorded as the
# translator in column D.
$rowsAddTranslator = @(378,379,380,381,382,383,384,385,386,387,388,389,390,391,392,393,394,395,396,397,398,399,400,401,427,428,429,430,431,432,434,435,436,437,439,440,441,442,443,445)

foreach ($r in $rowsAddTranslator) {
    $ws.Cells.Item($r, 4).Value = "juminho"
}

# Update the view state to match where the editor ended up working.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 425
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D445").Select()
